# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Update status text shared across sheets (Overview, zh-cn, de-de) from
# "Ready for handoff" to "Handback transform failed".
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handback transform failed"
$overview.Range("F3").Value = "Handback transform failed"
$zhcn.Range("C3").Value = "Handback transform failed"
$dede.Range("C3").Value = "Handback transform failed"

# Record the handback/handoff filename mismatch error detail for row 3
# (source file 7f7a080f-9052-493f-9804-faa480bc1a1e.md) on both locale sheets.
$zhcn.Range("P3").Value = "Handback file name: fixcirv5.myj is different with handoff file name: 7f7a080f-9052-493f-9804-faa480bc1a1e.3eba3c56d7cde081f24a9565af8a6005f9bebe38.zh-cn."
$dede.Range("P3").Value = "Handback file name: fixcirv5.myj is different with handoff file name: 7f7a080f-9052-493f-9804-faa480bc1a1e.3eba3c56d7cde081f24a9565af8a6005f9bebe38.de-de."

# Widen the Error Detail column to fit the new message text (target stored
# XML width of 40; the COM ColumnWidth setter adds the standard ~5/6 char
# padding term when converting to the stored grid width, so back it out).
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
